$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# The household_id row (row 2) becomes a read-only "note" field instead of
# an editable "text" field. It loses its element name (column C, "name")
# and gets a new display.text (column D) that echoes the value instead of
# prompting for it.
$ws.Range("A2").Value = "note"
$ws.Range("C2").ClearContents() | Out-Null
$ws.Range("D2").Value = "Data for household: {{household_id}}"

# The "readonly" and "comments" columns (F and G) are no longer used now
# that "readonly" is modeled via the note type itself, so drop both
# columns in one go (this shifts the old "hideInContents" column, H, left
# into F).
$ws.Range("F1:G1").EntireColumn.Delete() | Out-Null

# Leave the selection where the hideInContents flag for the new note row
# now lives.
$ws.Range("F6").Select() | Out-Null
